$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Edison -> IEC" column (column O) entirely. This shifts
# "24' Extension", "Power Strip" and "Tri-Tap" one column to the left
# (P->O, Q->P, R->Q) and drops the now-empty trailing column.
$ws.Range("O1").EntireColumn.Delete()

# Update the DMX cable length labels per Ray's feedback.
$ws.Range("G2").Value = "5’"
$ws.Range("H2").Value = "10’"
$ws.Range("I2").Value = "25’"
$ws.Range("J2").Value = "50’"

# Fix the DMX cable counts for "Totems x4".
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 3

# Fix the 3/5 Pin Converters / Edison -> True1 counts for "12P Hex":
# the 12 converters were actually Edison -> True1, not 3/5 Pin Converters.
$ws.Range("K7").ClearContents()
$ws.Range("N7").Value = 12

# Move the active cell selection to match the edited workbook.
$ws.Range("D14").Select()
